$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("G3").Value = 2
$ws.Range("H3").Value = 3.4
$ws.Range("I3").Value = 3.6
$ws.Range("J3").Value = 1.04
$ws.Range("K3").Value = 12
$ws.Range("L3").Value = 1.22
$ws.Range("M3").Value = 4
$ws.Range("N3").Value = 1.75
$ws.Range("O3").Value = 2.05
$ws.Range("P3").Value = 1.33
$ws.Range("Q3").Value = 3.25
$ws.Range("R3").Value = 1.62
$ws.Range("S3").Value = 2.2
$ws.Range("T3").Value = 9
$ws.Range("U3").Value = 11
$ws.Range("V3").Value = 9
$ws.Range("W3").Value = 19
$ws.Range("Y3").Value = 23
$ws.Range("Z3").Value = 12
$ws.Range("AA3").Value = 6.5
$ws.Range("AB3").Value = 12
$ws.Range("AC3").Value = 41
$ws.Range("AD3").Value = 151
$ws.Range("AE3").Value = 13
$ws.Range("AI3").Value = 26
$ws.Range("AJ3").Value = 34

# Row 5 updates
$ws.Range("H5").Value = 3.9
$ws.Range("N5").Value = 1.6
$ws.Range("V5").Value = 16.5
$ws.Range("W5").Value = 100
$ws.Range("X5").Value = 50
$ws.Range("Y5").Value = 45
$ws.Range("AA5").Value = 7.8
$ws.Range("AJ5").Value = 22
